$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9: Distill, My Heart / Distilled Water (item 5487)
$ws.Range("H9").Value = 150.6
$ws.Range("I9").Value = 157.91667
$ws.Range("J9").Value = 121.333336
$ws.Range("K9").Value = 157.91667
$ws.Range("L9").Value = 121.333336
$ws.Range("M9").Value = 11.08332999999999
$ws.Range("N9").Value = -459.333336

# Row 17: One for the Road / Potion (item 38956)
$ws.Range("H17").Value = 1479.5186
$ws.Range("J17").Value = 1543.875
$ws.Range("L17").Value = 4631.625
$ws.Range("N17").Value = -4967.625

# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink (item 27781)
$ws.Range("H62").Value = 21526.092
$ws.Range("I62").Value = 21420.111
$ws.Range("K62").Value = 21420.111
$ws.Range("M62").Value = -20796.111

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink (item 27781)
$ws.Range("H65").Value = 21526.092
$ws.Range("I65").Value = 21420.111
$ws.Range("K65").Value = 107100.555
$ws.Range("M65").Value = -103980.555

# Row 112: Making Ends Meet / Superior Spiritbond Potion (item 27960)
$ws.Range("H112").Value = 1151.3636
$ws.Range("J112").Value = 1173.2
$ws.Range("L112").Value = 3519.6
$ws.Range("N112").Value = -5735.6

# Row 129: Practical Command / Commanding Craftsman's Draught (item 36115)
$ws.Range("H129").Value = 2140.3125
$ws.Range("J129").Value = 3000
$ws.Range("L129").Value = 9000
$ws.Range("N129").Value = -19000

# Row 132: Fast-forwarding Flora / Growth Formula Lambda (item 44049)
$ws.Range("H132").Value = 41202.62
$ws.Range("I132").Value = 47539.746
$ws.Range("J132").Value = 8883.299999999999
$ws.Range("K132").Value = 142619.238
$ws.Range("L132").Value = 26649.9
$ws.Range("M132").Value = -140089.238
$ws.Range("N132").Value = -31709.9

# Row 133: Big Brush, Big Dreams / Ginseng Angle Brush (item 41856)
$ws.Range("H133").Value = 85706.28999999999
$ws.Range("J133").Value = 85706.28999999999
$ws.Range("L133").Value = 85706.28999999999
$ws.Range("N133").Value = -95826.28999999999

# Row 138: All-night Crafting / Cunning Craftsman's Tisane (item 44169)
$ws.Range("H138").Value = 2301.8235
$ws.Range("I138").Value = 1342.3125
$ws.Range("J138").Value = 2740.457
$ws.Range("K138").Value = 4026.9375
$ws.Range("L138").Value = 8221.370999999999
$ws.Range("M138").Value = 1113.0625
$ws.Range("N138").Value = -18501.371

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot (item 27713)
$ws.Range("H2").Value = 3435.4546
$ws.Range("I2").Value = 1977
$ws.Range("K2").Value = 1977
$ws.Range("M2").Value = -1864

# Row 14: Waste Not, Want Not / Bronze Celata (item 2673)
$ws.Range("H14").Value = 2103.2856
$ws.Range("I14").Value = 3657.5
$ws.Range("J14").Value = 937.625
$ws.Range("K14").Value = 3657.5
$ws.Range("L14").Value = 937.625
$ws.Range("M14").Value = -3482.5
$ws.Range("N14").Value = -1287.625

# Row 32: Ingot We Trust / Steel Ingot (item 44147)
$ws.Range("H32").Value = 6173647.5
$ws.Range("I32").Value = 6250816.5
$ws.Range("K32").Value = 6250816.5
$ws.Range("M32").Value = -6250529.5

# Row 45: Hollow Hallmarks / Mythril Ingot (item 27714)
$ws.Range("H45").Value = 2128.4
$ws.Range("I45").Value = 2073.739
$ws.Range("K45").Value = 2073.739
$ws.Range("M45").Value = -1696.739

# Row 74: As the Bolt Flies / Titanium Nugget (item 44000)
$ws.Range("H74").Value = 4033886.5
$ws.Range("I74").Value = 5000963
$ws.Range("J74").Value = 4399.8335
$ws.Range("K74").Value = 5000963
$ws.Range("L74").Value = 4399.8335
$ws.Range("M74").Value = -5000089
$ws.Range("N74").Value = -6147.8335

# Row 77: Heavy Metal Banned (L) / Titanium Nugget (item 44000)
$ws.Range("H77").Value = 4033886.5
$ws.Range("I77").Value = 5000963
$ws.Range("J77").Value = 4399.8335
$ws.Range("K77").Value = 25004815
$ws.Range("L77").Value = 21999.1675
$ws.Range("M77").Value = -25000447
$ws.Range("N77").Value = -30735.1675

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot (item 19945)
$ws.Range("H102").Value = 52677.57
$ws.Range("I102").Value = 60457.332
$ws.Range("K102").Value = 60457.332
$ws.Range("M102").Value = -58835.332

# Row 110: Scheduled Maintenance / Deepgold Ingot (item 27708)
$ws.Range("H110").Value = 1718.75
$ws.Range("I110").Value = 1691.6666
$ws.Range("K110").Value = 1691.6666
$ws.Range("M110").Value = 353.3334

# Row 116: No Scope / Titanbronze Ingot (item 27713)
$ws.Range("H116").Value = 3435.4546
$ws.Range("I116").Value = 1977
$ws.Range("K116").Value = 1977
$ws.Range("M116").Value = 317

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot (item 43997)
$ws.Range("H132").Value = 397366.2
$ws.Range("I132").Value = 569909.4
$ws.Range("J132").Value = 4351.1665
$ws.Range("K132").Value = 1709728.2
$ws.Range("L132").Value = 13053.4995
$ws.Range("M132").Value = -1707198.2
$ws.Range("N132").Value = -18113.4995

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot (item 27713)
$ws.Range("H3").Value = 3435.4546
$ws.Range("I3").Value = 1977
$ws.Range("K3").Value = 1977
$ws.Range("M3").Value = -1863

# Row 86: Through Thick and Thin / Adamantite Nugget (item 12526)
$ws.Range("H86").Value = 1644.6666
$ws.Range("I86").Value = 1483.1765
$ws.Range("J86").Value = 1855.8462
$ws.Range("K86").Value = 1483.1765
$ws.Range("L86").Value = 1855.8462
$ws.Range("M86").Value = -360.1765
$ws.Range("N86").Value = -4101.8462

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget (item 12526)
$ws.Range("H89").Value = 1644.6666
$ws.Range("I89").Value = 1483.1765
$ws.Range("J89").Value = 1855.8462
$ws.Range("K89").Value = 7415.8825
$ws.Range("L89").Value = 9279.231
$ws.Range("M89").Value = -1799.8825
$ws.Range("N89").Value = -20511.231

# Row 99: Meddle in Metal / Oroshigane Ingot (item 19943)
$ws.Range("H99").Value = 13207.333
$ws.Range("I99").Value = 6684.875
$ws.Range("K99").Value = 6684.875
$ws.Range("M99").Value = -5186.875

# Row 105: Ingot to Wing It / Molybdenum Ingot (item 19947)
$ws.Range("H105").Value = 2203.25
$ws.Range("I105").Value = 2104.3333
$ws.Range("K105").Value = 2104.3333
$ws.Range("M105").Value = -357.3332999999998

# Row 107: The Gold Experience / Deepgold Nugget (item 27706)
$ws.Range("H107").Value = 6518.3335
$ws.Range("J107").Value = 4500
$ws.Range("L107").Value = 4500
$ws.Range("N107").Value = -8340

$ws = $wb.Worksheets.Item("CRP")
# Row 86: Birch, Please / Birch Lumber (item 12584)
$ws.Range("H86").Value = 340661.34
$ws.Range("I86").Value = 8491.5
$ws.Range("K86").Value = 8491.5
$ws.Range("M86").Value = -7368.5

# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber (item 12584)
$ws.Range("H89").Value = 340661.34
$ws.Range("I89").Value = 8491.5
$ws.Range("K89").Value = 42457.5
$ws.Range("M89").Value = -36841.5

# Row 105: Zelkova, My Love / Zelkova Lumber (item 19928)
$ws.Range("H105").Value = 33303.09
$ws.Range("I105").Value = 44869.375
$ws.Range("J105").Value = 2459.6667
$ws.Range("K105").Value = 44869.375
$ws.Range("L105").Value = 2459.6667
$ws.Range("M105").Value = -43122.375
$ws.Range("N105").Value = -5953.6667

# Row 107: Built to Last / White Oak Lumber (item 27689)
$ws.Range("H107").Value = 1128.1666
$ws.Range("I107").Value = 1017.75
$ws.Range("J107").Value = 1349
$ws.Range("K107").Value = 1017.75
$ws.Range("L107").Value = 1349
$ws.Range("M107").Value = 902.25
$ws.Range("N107").Value = -5189

# Row 132: Hull Lotta Damage / Ginseng Lumber (item 44019)
$ws.Range("H132").Value = 86557064
$ws.Range("I132").Value = 142885900
$ws.Range("K132").Value = 428657700
$ws.Range("M132").Value = -428655170

$ws = $wb.Worksheets.Item("CUL")
# Row 46: Feeding Frenzy / Acorn Cookie (item 4701)
$ws.Range("H46").Value = 3
$ws.Range("I46").Value = 2.5
$ws.Range("J46").Value = 4
$ws.Range("K46").Value = 7.5
$ws.Range("L46").Value = 12
$ws.Range("M46").Value = 83.5
$ws.Range("N46").Value = -194

# Row 122: Salt of the North / Northern Sea Salt (item 36078)
$ws.Range("H122").Value = 20000194
$ws.Range("I122").Value = 199.42857
$ws.Range("J122").Value = 66666850
$ws.Range("K122").Value = 1794.85713
$ws.Range("L122").Value = 600001650
$ws.Range("M122").Value = 655.1428699999999
$ws.Range("N122").Value = -600006550

# Row 124: Bobbing for Compliments / Island Miq'abob (item 36040)
$ws.Range("H124").Value = 1236.5
$ws.Range("I124").Value = 1236.5
$ws.Range("K124").Value = 3709.5
$ws.Range("M124").Value = 1200.5

# Row 125: At Any Temperature / Borscht (item 36043)
$ws.Range("H125").Value = 17808.857
$ws.Range("I125").Value = 6015
$ws.Range("J125").Value = 22526.4
$ws.Range("K125").Value = 18045
$ws.Range("L125").Value = 67579.20000000001
$ws.Range("M125").Value = -13125
$ws.Range("N125").Value = -77419.20000000001

# Row 126: Imperial Palate / Glory Be Soup (item 36045)
$ws.Range("H126").Value = 3419.2
$ws.Range("I126").Value = 4010
$ws.Range("J126").Value = 2533
$ws.Range("K126").Value = 12030
$ws.Range("L126").Value = 7599
$ws.Range("M126").Value = -7090
$ws.Range("N126").Value = -17479

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle / Durium Ingot (item 36169)
$ws.Range("H102").Value = 1963.4667
$ws.Range("I102").Value = 1845.3
$ws.Range("J102").Value = 2199.8
$ws.Range("K102").Value = 1845.3
$ws.Range("L102").Value = 2199.8
$ws.Range("M102").Value = -223.3
$ws.Range("N102").Value = -5443.8

# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone (item 27802)
$ws.Range("H107").Value = 52421.855
$ws.Range("J107").Value = 3802
$ws.Range("L107").Value = 3802
$ws.Range("N107").Value = -7642

# Row 113: Copious Crystal Cannons / Manasilver Nugget (item 27710)
$ws.Range("H113").Value = 3750
$ws.Range("J113").Value = 6000
$ws.Range("L113").Value = 6000
$ws.Range("N113").Value = -10340

# Row 123: Workplace Workout / Ametrine Ring of Fending (item 34150)
$ws.Range("H123").Value = 61417
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 61417
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 61417
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -66317

# Row 132: On Board for Lar / Lar Ingot (item 44008)
$ws.Range("H132").Value = 16075526
$ws.Range("I132").Value = 24696314
$ws.Range("K132").Value = 74088942
$ws.Range("M132").Value = -74086412

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather (item 5277)
$ws.Range("H22").Value = 552.1429000000001
$ws.Range("I22").Value = 480.2
$ws.Range("J22").Value = 732
$ws.Range("K22").Value = 480.2
$ws.Range("L22").Value = 732
$ws.Range("M22").Value = -185.2
$ws.Range("N22").Value = -1322

# Row 27: Fire and Hide / Aldgoat Leather (item 5277)
$ws.Range("H27").Value = 552.1429000000001
$ws.Range("I27").Value = 480.2
$ws.Range("J27").Value = 732
$ws.Range("K27").Value = 480.2
$ws.Range("L27").Value = 732
$ws.Range("M27").Value = -373.2
$ws.Range("N27").Value = -946

# Row 74: Overall, We Blend In / Dhalmelskin Vest (item 11990)
$ws.Range("H74").Value = 77173.60000000001
$ws.Range("J74").Value = 77173.60000000001
$ws.Range("L74").Value = 77173.60000000001
$ws.Range("N74").Value = -79169.60000000001

# Row 77: Eviction Notice (L) / Dhalmelskin Vest (item 11990)
$ws.Range("H77").Value = 77173.60000000001
$ws.Range("J77").Value = 77173.60000000001
$ws.Range("L77").Value = 231520.8
$ws.Range("N77").Value = -241504.8

# Row 122: Hell on Leather / Gaja Leather (item 36247)
$ws.Range("H122").Value = 5030.6665
$ws.Range("I122").Value = 4596.722
$ws.Range("J122").Value = 6332.5
$ws.Range("K122").Value = 13790.166
$ws.Range("L122").Value = 18997.5
$ws.Range("M122").Value = -11340.166
$ws.Range("N122").Value = -23897.5

$ws = $wb.Worksheets.Item("WVR")
# Row 100: Of Great Import / Kudzu Thread (item 19981)
$ws.Range("H100").Value = 1534.68
$ws.Range("I100").Value = 1072.0667
$ws.Range("K100").Value = 2144.1334
$ws.Range("M100").Value = -1603.1334
